$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Headers (row 1) ----
$ws.Range("A1").Value = "fecha"
$ws.Range("B1").Value = "camaronera"
$ws.Range("C1").Value = "clase"
$ws.Range("D1").Value = "recibidas"
$ws.Range("E1").Value = "reparadas"
$ws.Range("F1").Value = "descartadas"
$ws.Range("G1").Value = "responsable"

# New header cells (F1:G1) need the same look as the existing bold/bordered
# header style already applied to A1:E1 - copy formats from an existing
# header cell so no new font/style records are introduced.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Data rows (2-9) ----
# fecha, camaronera, clase, recibidas, reparadas, descartadas, responsable
$data = @(
    @("2024-01-24", "AFRICA",   "AMA", 2,  2, 2,  "JORGE FUENTES"),
    @("2024-01-24", "AFRICA ",  "ASP", 30, 2, 3,  "JORGE FUENTES"),
    @("2024-01-24", "AFRICA",   "ASP", 2,  2, 2,  "JORGE FUENTES"),
    @("2024-01-24", "3",        "ASP", 2,  1, 0,  "JORGE FUENTES"),
    @("2024-01-24", "BAJEN",    "ASP", 30, 2, 20, "JORGE FUENTES"),
    @("2024-01-24", "AFRICA",   "ASP", 2,  2, 2,  "JORGE FUENTES"),
    @("2024-01-24", "AFRICA ",  "ASP", 2,  2, 5,  "JORGE FUENTES"),
    @("2024-01-24", "AFRICA",   "ASP", 2,  2, 2,  "JORGE FUENTES")
)

$row = 2
foreach ($r in $data) {
    # Column A ("fecha") holds text that looks like a date; force it to stay
    # text (not get auto-converted to a date serial number) while keeping
    # the cell's style at the default (no explicit number format applied).
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $r[0]
    $cellA.ClearFormats()

    # Column B ("camaronera") is usually plain text, but row 5 holds the
    # numeric-looking "3" which must also stay text.
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $r[1]
    $cellB.ClearFormats()

    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row++
}
